$d = $word.ActiveDocument

# The "Requisitos" list paragraph currently lists:
#   LOQ4095 -  Química Geral Experimental  (Requisito fraco)
#   LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
#   LOB1012 -  Estatística  (Requisito fraco)
# Each line is its own run ending with a manual line break (<w:br/>).
# We need to move the "LOB1012" line from the end of the list to the front,
# immediately before the "LOQ4095" line, without disturbing the other lines.

$loq4095 = "LOQ4095 -  Qu" + [char]0x00ED + "mica Geral Experimental  (Requisito fraco)"
$lob1012 = "LOB1012 -  Estat" + [char]0x00ED + "stica  (Requisito fraco)"
$lineBreak = [char]11

# 1) Find the start of the "LOQ4095" line -- that's where the moved line
#    needs to be inserted (as a new run, so it stays distinct from the
#    following run).
$findRange = $d.Content
$found = $findRange.Find.Execute($loq4095, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the LOQ4095 requirement line."
}
$insertAt = $findRange.Start

# 2) Insert the LOB1012 line (with its own trailing manual line break) right
#    before the LOQ4095 line.
$insertionPoint = $d.Range($insertAt, $insertAt)
$insertionPoint.InsertBefore($lob1012 + $lineBreak)

# 3) Remove the original LOB1012 line (and its trailing line break) further
#    down in the same paragraph, now shifted forward by the length of what
#    was just inserted.
$searchFrom = $insertAt + $lob1012.Length + 1
$tailRange = $d.Range($searchFrom, $d.Content.End)
$found2 = $tailRange.Find.Execute($lob1012, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the original LOB1012 requirement line to remove."
}
$deleteRange = $d.Range($tailRange.Start, $tailRange.End + 1)
$deleteRange.Delete()
